$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy rows 23:24 (the "test case 5" block) down to rows 27:28 to create "test case 6"
$ws.Range("A23:G24").Copy($ws.Range("A27"))

# Update the text of the newly pasted cells for the new "test case 6" scenario.
# Shared-string insertion order matters, so write in this exact order:
# A27, G27, C27, C28
$ws.Range("A27").Value = "test case 6"
$ws.Range("G27").Value = "Method int test6(int a, String b, Double[] c)"
$ws.Range("C27").Value = "Method int test6Arguments()"
$ws.Range("C28").Value = "return test6(5, ""abc"", 10, (long) 12, (double) 14);"

# Update selection to match the final state
$ws.Range("C30").Select()
